# Apply "Updated plans for Marketing" edits to the School Visit Plan sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily-Marketing-Plan")

# New school entries appended after the existing 20 rows (row 23 onward).
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "Calcutta Public School"
$ws.Cells.Item(24, 2).Value = "St. Ignatius"
$ws.Cells.Item(25, 2).Value = "Young Horizon"
$ws.Cells.Item(26, 2).Value = "Swarnim "
$ws.Cells.Item(27, 2).Value = "Lion's "
$ws.Cells.Item(28, 2).Value = "John Bosco"
$ws.Cells.Item(29, 2).Value = "St Stephens"
$ws.Cells.Item(30, 2).Value = "AG Tollygunj"

# Clear the stray formatting on G5 (it no longer carries a distinct style).
$ws.Range("G5").Style = "Normal"

$wb.Save()
